$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row: Right marks per question 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row: Right total 60 -> 100
$ws.Range("B12").Value = 100

# Update "Total" row: Correct/Total marks display 56/84 -> 100/140
$ws.Range("E12").Value = "100/140"
